$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlShiftToLeft = -4159

# --- Header row updates ---
$ws.Range("H1").Value = "In Hospital Now"
$ws.Range("I1").Value = "Total Been In Hospital"

# Remove the now-unused trailing header columns (Community Transmission,
# Total Community Transmission, Tests, Total Tests, Test Count) and shift
# the (empty) cells beyond them left so the row - and sheet dimension -
# truly ends at column J.
$ws.Range("K1:O1").Delete($xlShiftToLeft)

# --- Row 26: drop Community Transmission / Tests columns, keep A:C only ---
$ws.Range("K26").Value = $null
$ws.Range("L26").Value = $null
$ws.Range("N26").Value = $null
$ws.Range("O26").Value = $null

# --- Row 27: move "In Hospital" value from I27 to H27, drop I,K:O ---
$ws.Range("H27").Value = $ws.Range("I27").Value2
$ws.Range("I27").Value = $null
$ws.Range("K27").Value = $null
$ws.Range("L27").Value = $null
$ws.Range("M27").Value = $null
$ws.Range("N27").Value = $null
$ws.Range("O27").Value = $null

# --- Row 28: move "In Hospital" value from I28 to H28, drop I,K:O ---
$ws.Range("H28").Value = $ws.Range("I28").Value2
$ws.Range("I28").Value = $null
$ws.Range("K28").Value = $null
$ws.Range("L28").Value = $null
$ws.Range("M28").Value = $null
$ws.Range("N28").Value = $null
$ws.Range("O28").Value = $null

# --- Row 29: replace H29 with old I29 value, drop I,K:O ---
$ws.Range("H29").Value = $ws.Range("I29").Value2
$ws.Range("I29").Value = $null
$ws.Range("K29").Value = $null
$ws.Range("L29").Value = $null
$ws.Range("M29").Value = $null
$ws.Range("N29").Value = $null
$ws.Range("O29").Value = $null

# --- New row 30: 27 March figures ---
$ws.Range("A30").Value = 43917
$ws.Range("B30").Value = 76
$ws.Range("C30").Value = 338
$ws.Range("D30").Value = 9
$ws.Range("E30").Value = 30
$ws.Range("F30").Value = 10
$ws.Range("G30").Value = 37
$ws.Range("H30").Value = 8
$ws.Range("I30").Value = 20
$ws.Range("J30").Value = 1

# Match the date style/format used by the other date cells in column A
$ws.Range("A30").NumberFormat = $ws.Range("A29").NumberFormat

# --- Dimension / used range should now be A1:J30 ---
